$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet / "Through date" label
$ws.Name = "Through 2022-05-02"

# Update the header label for the "through" year-to-date column (I1), which
# references the shared string "2022 (through 05-01)"
$ws.Range("I1").Value = "2022 (through 05-02)"

# Update the June row's year-to-date value (new May data point added)
$ws.Range("I6").Value = 4

# Update the Total row's year-to-date value
$ws.Range("I14").Value = 555
